$wb = $excel.ActiveWorkbook

# Sheet "zh-cn": update Correspond Handoff/Handback Datetime for row 2
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E2").Value = "2016-03-22 19:11:28"
$wsZh.Range("H2").Value = "2016-03-22 19:11:49"

# Sheet "de-de": update Correspond Handoff/Handback Datetime for row 2
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E2").Value = "2016-03-22 19:11:31"
$wsDe.Range("H2").Value = "2016-03-22 19:11:57"
